# Update "paises.xlsx" COVID data snapshot (3 Oct 2020, 15:58 -> 17:15):
#  - refresh the "last updated" timestamp
#  - refresh totals/new-cases/active/recovered/critical/deaths for the
#    countries whose stats moved between snapshots
#  - a handful of countries swapped rank/row order as their totals shifted,
#    so the country name in column A is rewritten for those rows too
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 17:15"
$ws.Range("B4").Value = 7557789
$ws.Range("C4").Value = 8466
$ws.Range("D4").Value = 4777871
$ws.Range("E4").Value = 2566323
$ws.Range("G4").Value = 71
$ws.Range("H4").Value = 213595

$ws.Range("B5").Value = 6509916
$ws.Range("C5").Value = 37982
$ws.Range("D5").Value = 5466344
$ws.Range("E5").Value = 942361
$ws.Range("G5").Value = 336
$ws.Range("H5").Value = 101211

$ws.Range("A15").Value = "Chile"
$ws.Range("B15").Value = 468471
$ws.Range("C15").Value = 1881
$ws.Range("D15").Value = 440881
$ws.Range("E15").Value = 14671
$ws.Range("G15").Value = 52
$ws.Range("H15").Value = 12919

$ws.Range("A16").Value = "Iran"
$ws.Range("B16").Value = 468119
$ws.Range("C16").Value = 3523
$ws.Range("D16").Value = 387675
$ws.Range("E16").Value = 53698
$ws.Range("G16").Value = 179
$ws.Range("H16").Value = 26746

$ws.Range("A17").Value = "Reino Unido"
$ws.Range("B17").Value = 467146
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("H17").Value = 42268

$ws.Range("A21").Value = "Italia"
$ws.Range("B21").Value = 322751
$ws.Range("C21").Value = 2844
$ws.Range("D21").Value = 231217
$ws.Range("E21").Value = 55566
$ws.Range("G21").Value = 27
$ws.Range("H21").Value = 35968

$ws.Range("A22").Value = "Turquia"
$ws.Range("B22").Value = 321512
$ws.Range("D22").Value = 282657
$ws.Range("E22").Value = 30530
$ws.Range("H22").Value = 8325

$ws.Range("B26").Value = 298475
$ws.Range("C26").Value = 112
$ws.Range("E26").Value = 29379

$ws.Range("B38").Value = 113926
$ws.Range("C38").Value = 576
$ws.Range("D38").Value = 89635
$ws.Range("E38").Value = 22163
$ws.Range("G38").Value = 11
$ws.Range("H38").Value = 2128

$ws.Range("B46").Value = 93748
$ws.Range("C46").Value = 658
$ws.Range("D46").Value = 82040
$ws.Range("E46").Value = 8423
$ws.Range("G46").Value = 18
$ws.Range("H46").Value = 3285

$ws.Range("A48").Value = "Japon"
$ws.Range("B48").Value = 84768
$ws.Range("C48").Value = 553
$ws.Range("D48").Value = 77807
$ws.Range("E48").Value = 5371
$ws.Range("G48").Value = 12
$ws.Range("H48").Value = 1590

$ws.Range("A49").Value = "Nepal"
$ws.Range("B49").Value = 84570
$ws.Range("C49").Value = 2120
$ws.Range("D49").Value = 62740
$ws.Range("E49").Value = 21302
$ws.Range("G49").Value = 8
$ws.Range("H49").Value = 528

$ws.Range("A52").Value = "Portugal"
$ws.Range("B52").Value = 78247
$ws.Range("C52").Value = 963
$ws.Range("D52").Value = 49845
$ws.Range("E52").Value = 26407
$ws.Range("G52").Value = 12
$ws.Range("H52").Value = 1995

$ws.Range("A53").Value = "Chequia"
$ws.Range("B53").Value = 78051
$ws.Range("D53").Value = 35032
$ws.Range("E53").Value = 42320
$ws.Range("H53").Value = 699

$ws.Range("A54").Value = "Costa Rica"
$ws.Range("B54").Value = 77829
$ws.Range("D54").Value = 42621
$ws.Range("E54").Value = 34278
$ws.Range("H54").Value = 930

$ws.Range("A59").Value = "Uzbekistan"
$ws.Range("B59").Value = 57924
$ws.Range("C59").Value = 470
$ws.Range("D59").Value = 54456
$ws.Range("E59").Value = 2992
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 476

$ws.Range("A60").Value = "Singapur"
$ws.Range("B60").Value = 57800
$ws.Range("C60").Value = 6
$ws.Range("D60").Value = 57562
$ws.Range("E60").Value = 211
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 27

$ws.Range("B61").Value = 55888
$ws.Range("C61").Value = 872
$ws.Range("D61").Value = 40475
$ws.Range("E61").Value = 14060
$ws.Range("G61").Value = 9
$ws.Range("H61").Value = 1353

$ws.Range("B95").Value = 14117
$ws.Range("C95").Value = 152
$ws.Range("D95").Value = 8536
$ws.Range("E95").Value = 5189
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 392

$ws.Range("B120").Value = 5780
$ws.Range("C120").Value = 62
$ws.Range("D120").Value = 5118
$ws.Range("E120").Value = 540

$ws.Range("B174").Value = 540
$ws.Range("C174").Value = 1
$ws.Range("D174").Value = 527
$ws.Range("E174").Value = 6

$ws.Range("A207").Value = "Nueva Caledonia"

$ws.Range("A208").Value = "Santa Lucia"

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

